# Fix bug in convertJsonToExcel: refresh "last_edited_time" (column D) for
# the Notion rows that were re-synced, and correct the Doanh so (AO) /
# Ti le dat KPI (BI) figures that the buggy export had wrong for 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-26 all now share the same last_edited_time stamp (20:33), row 27
# keeps the previous minute's stamp (20:32) but re-synced to the new day.
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-08-24T20:33:00.000Z"
}
$ws.Cells.Item(27, 4).Value = "2024-08-24T20:32:00.000Z"

# Corrected sales (Doanh so) and KPI ratio (Ti le dat KPI) figures.
$ws.Range("AO21").Value = 41800000
$ws.Range("BI21").Value = 1.3933

$ws.Range("AO22").Value = 106500000
$ws.Range("BI22").Value = 3.55

$ws.Range("AO25").Value = 118500000
$ws.Range("BI25").Value = 3.95
